$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 150 (pushes old row 150 down to row 152)
$ws.Rows("150:151").Insert()

# --- Update existing rows 146-149 with their new values ---

# Row 146
$ws.Cells.Item(146, 4).Value = 44595
$ws.Cells.Item(146, 10).Value = 600
$ws.Cells.Item(146, 11).Value = 28000
$ws.Cells.Item(146, 12).Value = 31000
$ws.Cells.Item(146, 13).Value = 29500
$ws.Cells.Item(146, 16).Value = 1180

# Row 147
$ws.Cells.Item(147, 4).Value = 44595
$ws.Cells.Item(147, 10).Value = 300
$ws.Cells.Item(147, 11).Value = 35000
$ws.Cells.Item(147, 12).Value = 40000
$ws.Cells.Item(147, 13).Value = 37500
$ws.Cells.Item(147, 16).Value = 1500

# Row 148
$ws.Cells.Item(148, 4).Value = 44335
$ws.Cells.Item(148, 10).Value = 1000
$ws.Cells.Item(148, 11).Value = 26000
$ws.Cells.Item(148, 12).Value = 27000
$ws.Cells.Item(148, 13).Value = 26500
$ws.Cells.Item(148, 16).Value = 1060

# Row 149
$ws.Cells.Item(149, 4).Value = 44335
$ws.Cells.Item(149, 10).Value = 800
$ws.Cells.Item(149, 11).Value = 30000
$ws.Cells.Item(149, 12).Value = 31000
$ws.Cells.Item(149, 13).Value = 30500
$ws.Cells.Item(149, 16).Value = 1220

# --- Fill the two newly inserted rows (150, 151) ---

# Row 150
$ws.Cells.Item(150, 1).Value = 2
$ws.Cells.Item(150, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44552
$ws.Cells.Item(150, 5).Value = 4
$ws.Cells.Item(150, 6).Value = 100112031
$ws.Cells.Item(150, 7).Value = "Poroto verde"
$ws.Cells.Item(150, 8).Value = "Magnum"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 600
$ws.Cells.Item(150, 11).Value = 15000
$ws.Cells.Item(150, 12).Value = 16000
$ws.Cells.Item(150, 13).Value = 15500
$ws.Cells.Item(150, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(150, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(150, 16).Value = 620
$ws.Cells.Item(150, 17).Value = 25
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# Row 151
$ws.Cells.Item(151, 1).Value = 2
$ws.Cells.Item(151, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44552
$ws.Cells.Item(151, 5).Value = 4
$ws.Cells.Item(151, 6).Value = 100112031
$ws.Cells.Item(151, 7).Value = "Poroto verde"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 500
$ws.Cells.Item(151, 11).Value = 28000
$ws.Cells.Item(151, 12).Value = 30000
$ws.Cells.Item(151, 13).Value = 29000
$ws.Cells.Item(151, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(151, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(151, 16).Value = 1160
$ws.Cells.Item(151, 17).Value = 25
$ws.Cells.Item(151, 18).Value = "Hortaliza"
